$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" strings look like plain decimal numbers (e.g.
# "594.05"), unlike most of the source data which uses thousand-separator dots
# (e.g. "3.656.06") and is therefore never auto-converted. Excel would silently
# reinterpret those plain-looking numeric strings as numbers on assignment, so
# each such cell is pre-formatted as Text to preserve the literal string value.
$textCells = @("D5", "D6", "D11", "D12", "D13", "D18", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D32", "D33", "D37", "D38", "D43", "D44", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.169.48"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").Value = "3.656.06"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "594.05"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").Value = "165.29"
$ws.Range("E6").Value = "  -6.88%  "
$ws.Range("D7").Value = "3.649.74"
$ws.Range("E7").Value = "  -4.76%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -5.47%  "
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -5.03%  "
$ws.Range("D13").Value = "37.33"
$ws.Range("E13").Value = "  -6.28%  "
$ws.Range("E14").Value = "  -6.58%  "
$ws.Range("D15").Value = "4.274.01"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("D16").Value = "3.660.91"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "67.213.99"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  -5.68%  "
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "17.23"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "490.52"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").Value = "9.06"
$ws.Range("E22").Value = "  -6.13%  "
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").Value = "85.25"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("E25").Value = "  -7.48%  "
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  -5.95%  "
$ws.Range("D27").Value = "12.11"
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "9.94"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("E31").Value = "  -6.62%  "
$ws.Range("D32").Value = "7.64"
$ws.Range("D33").Value = "31.40"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "3.801.50"
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -6.36%  "
$ws.Range("D36").Value = "3.598.03"
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.986"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("E40").Value = "  -6.96%  "
$ws.Range("E41").Value = "  -3.96%  "
$ws.Range("E42").Value = "  -10.49%  "
$ws.Range("D43").Value = "48.61"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").Value = "1.91"
$ws.Range("E44").Value = "  -7.21%  "
$ws.Range("E45").Value = "  -9.48%  "
$ws.Range("D46").Value = "8.30"
$ws.Range("E46").Value = "  -3.00%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "40.47"
$ws.Range("E48").Value = "  -7.84%  "
$ws.Range("D49").Value = "142.26"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "2.743.37"
$ws.Range("E50").Value = "  -6.50%  "
$ws.Range("D51").Value = "0.0346"
$ws.Range("E51").Value = "  -4.41%  "
